$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.242.18"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "2.268.25"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "305.99"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "97.23"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  -0.98%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").Value = "35.16"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("E11").Value = "  -2.65%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "6.93"
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("D14").Value = "2.621.01"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "14.80"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "2.260.17"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "42.111.44"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "12.29"
$ws.Range("E19").Value = "  -4.91%  "
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").Value = "6.03"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "67.92"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "237.96"
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("E24").Value = "  -1.99%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "23.56"
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("D28").Value = "37.85"
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("D29").Value = "9.55"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").Value = "162.00"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").Value = "5.25"
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").Value = "17.69"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").Value = "2.37"
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("E38").Value = "  -4.59%  "
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("E41").Value = "  -3.88%  "
$ws.Range("E42").Value = "  +3.02%  "
$ws.Range("D43").Value = "19.41"
$ws.Range("E43").Value = "  -3.72%  "
$ws.Range("D44").Value = "1.948.99"
$ws.Range("E44").Value = "  -3.50%  "
$ws.Range("D45").Value = "0.0281"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("E46").Value = "  -3.04%  "
$ws.Range("D47").Value = "9.87"
$ws.Range("E47").Value = "  -3.99%  "
$ws.Range("D48").Value = "53.61"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").Value = "92.56"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").Value = "71.79"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  -2.64%  "
